$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply cell formatting for the two new 4-row test blocks (rows 70-73 and 74-77) ---
# Rows 57-60 already carry the target style pattern (s="4" on A/B/C/E, s="2" on D),
# so copy formats from there onto the new blocks.
$ws.Range("A57:E60").Copy()
$ws.Range("A70").PasteSpecial(-4122)
$ws.Range("A57:E60").Copy()
$ws.Range("A74").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 70-73: T_Dead_1 / Event.Deadline.Deadline ---
# Fill column-by-column (A, B, C, then D top-to-bottom, then E) to mirror the
# original authoring order reflected in the shared-string table.
$ws.Range("A70").Value = "T_Dead_1"
$ws.Range("B70").Value = "Event.Deadline.Deadline"
$ws.Range("C70").Value = "Verify user can create a deadline event with a deadline member."
$ws.Range("D70").Value = "1. Open Event creation form"
$ws.Range("D71").Value = "2. Select Deadline Event"
$ws.Range("D72").Value = "3. Input valid deadline event parameters and submit"
$ws.Range("D73").Value = "4. Verify the that the Deadline event object was created with the correct deadline member"
$ws.Range("E70").Value = "Deadline object is created with user given deadline member."

# --- Row 74-77: T_Dead_2 / Event.Deadline.WorkTime ---
$ws.Range("A74").Value = "T_Dead_2"
$ws.Range("B74").Value = "Event.Deadline.WorkTime"
$ws.Range("C74").Value = "Verify user can create a deadline event with a work time member."
$ws.Range("D74").Value = "1. Open Event creation form"
$ws.Range("D75").Value = "2. Select Deadline Event"
$ws.Range("D76").Value = "3. Input valid deadline event parameters and submit"
$ws.Range("D77").Value = "4. Verify the that the Deadline event object was created with the correct work time member"
$ws.Range("E74").Value = "Deadline object is created with user given work time member."

# --- Row heights matching the wrapped-text content ---
$ws.Rows.Item(70).RowHeight = 30
$ws.Rows.Item(72).RowHeight = 30
$ws.Rows.Item(73).RowHeight = 45
$ws.Rows.Item(76).RowHeight = 30
$ws.Rows.Item(77).RowHeight = 45

# --- Merge the per-test columns (A/B/C/E) across each 4-row block ---
$ws.Range("A70:A73").Merge()
$ws.Range("B70:B73").Merge()
$ws.Range("C70:C73").Merge()
$ws.Range("E70:E73").Merge()

$ws.Range("A74:A77").Merge()
$ws.Range("B74:B77").Merge()
$ws.Range("C74:C77").Merge()
$ws.Range("E74:E77").Merge()

# --- Update the visible selection to match the edited area ---
$ws.Range("C74:C77").Select()
